$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

# Title change (appears twice: Heading1 at top, and bold run near the end)
Replace-Text "Play Ancient Script Free - Review of Egyptian Slot Game" "Play Ancient Script Slot Game for Free"

# "What we like" bullet list
Replace-Text "Free spins and lucky symbol feature" "Egyptian theme with traditional graphics"
Replace-Text "Traditional graphics and immersive theme" "Special features for increased odds of success"
Replace-Text "Random symbol appearance adds excitement" "Immersive and exciting gameplay"
Replace-Text "Eye of Horus symbol offers big win potential" "Random symbol appearances for added excitement"

# "What we don't like" bullet list
Replace-Text "Big wins can be difficult to achieve" "Big wins may be elusive"

# Meta description italic paragraph near the end
Replace-Text "Explore the Ancient Script slot game, free to play, featuring an Egyptian theme with lucky symbols and free spins. Discover the pros and cons in our review." "Read our review of Ancient Script and play for free. Experience the Egyptian themed slot game with special features."
